# Team05Report.xlsx — "Added missing info for sprint 1 and LOC"
#
# 1. Team sheet: fill in Ian Klein's missing email address (D3) as a
#    mailto hyperlink, and make Team the active/selected sheet.
# 2. Burndown sheet: correct the LOC value for the second burndown entry
#    (D3), which recalculates the derived "code velocity" formula in F3.
# 3. Sprint1 sheet: fill in the missing actual size/time numbers for the
#    first three stories, and add the Sprint Review "Keep doing" / "Avoid"
#    notes (B22 / B26).

$wb = $excel.ActiveWorkbook

# --- Team sheet ---------------------------------------------------------
$team = $wb.Worksheets.Item("Team")
$team.Range("D3").Value = "iklein@stevens.edu"
$team.Hyperlinks.Add($team.Range("D3"), "mailto:iklein@stevens.edu")
$team.Activate()
$team.Range("D3").Select()

# --- Burndown sheet ------------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Range("D3").Value = 111
$burndown.Range("G6").Select()

# --- Sprint1 sheet ---------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

# Task JC01 actuals
$sprint1.Range("E2").Value = 20
$sprint1.Range("F2").Value = 10
$sprint1.Range("H2").Value = 20

# Task JC02 actuals
$sprint1.Range("E7").Value = 20
$sprint1.Range("F7").Value = 10
$sprint1.Range("H7").Value = 20

# Task JC03 actuals
$sprint1.Range("E12").Value = 20
$sprint1.Range("F12").Value = 20
$sprint1.Range("G12").Value = 19
$sprint1.Range("H12").Value = 20

# Sprint review notes
$sprint1.Range("B22").Value = "what we are doing.  Pretty consitent"
$sprint1.Range("B26").Value = "Making dumb mistakes.  Always error check when you go and there will be less bugs."

$sprint1.Range("D22").Select()

$team.Activate()
